# This script appends newly-observed daily snapshot rows to each of the
# three portfolio sheets in the workbook, mirroring the automated daily
# export job ("feat: 自动添加每日生成文件 AIPEEarningYield.xlsx (Colab)").
#
# For every sheet we append one row per holding (stock code/name/weight/
# quantity/price/value/total-value) plus the observation timestamp, using
# the same shape as all pre-existing rows.
#
# Columns: A=portfolio name, B=stock code, C=stock name, D=weight(%),
#          E=quantity, F=price, G=value, H=total value, I=timestamp.
#
# Numeric-looking text (stock codes such as "000333"/"100000" and the
# timestamp "202506121600") must be forced to Text, otherwise Excel's
# normal type inference would convert them into numbers (dropping
# leading zeros). We do this by temporarily setting NumberFormat="@"
# on those columns before writing the values, then calling
# ClearFormats() afterward so the cells end up with the same (default)
# style as every other data cell in the sheet.

function Add-PortfolioRows {
    param(
        $ws,
        [string]$label,
        [int]$startRow,
        [string]$timestamp,
        $rows
    )

    $n = $rows.Count
    $endRow = $startRow + $n - 1

    $fullRange = $ws.Range("A${startRow}:I${endRow}")
    $codeRange = $ws.Range("B${startRow}:B${endRow}")
    $timeRange = $ws.Range("I${startRow}:I${endRow}")

    # Force text storage for numeric-looking stock codes & timestamps.
    $codeRange.NumberFormat = "@"
    $timeRange.NumberFormat = "@"

    $arr = New-Object 'object[,]' $n,9
    for ($i = 0; $i -lt $n; $i++) {
        $r = $rows[$i]
        $arr[$i,0] = $label
        $arr[$i,1] = $r[0]
        $arr[$i,2] = $r[1]
        $arr[$i,3] = $r[2]
        $arr[$i,4] = $r[3]
        $arr[$i,5] = $r[4]
        $arr[$i,6] = $r[5]
        $arr[$i,7] = $r[6]
        $arr[$i,8] = $timestamp
    }
    $fullRange.Value = $arr

    # Drop back to default (unstyled) formatting, matching the rest of
    # the sheet's plain data rows.
    $codeRange.ClearFormats()
    $timeRange.ClearFormats()
}

$wb = $excel.ActiveWorkbook
$timestamp = "202506121600"


# Sheet 1: 大智投资组合收益 (新增第41-50行)
$ws = $wb.Worksheets.Item(1)
$rows = @(
    @("000333", "美的集团", 2.81, 39.71497483560155, 72.2, 2867.421183130432, 102146.0826454179),
    @("510050", "上证50ETF", 4.89, 1808.278443601665, 2.76, 4990.848504340595, 102146.0826454179),
    @("510300", "沪深300ETF", 4.9, 1247.712126085149, 4.01, 5003.325625601447, 102146.0826454179),
    @("518880", "黄金ETF", 4.95, 674.4389870730533, 7.5, 5058.2924030479, 102146.0826454179),
    @("600085", "同仁堂", 1.95, 52.96734947562319, 37.53, 1987.864625820138, 102146.0826454179),
    @("600900", "长江电力", 19.68, 665.8903941748626, 30.19, 20103.23100013911, 102146.0826454179),
    @("600989", "宝丰能源", 4.83, 306.7515982999751, 16.08, 4932.565700663599, 102146.0826454179),
    @("HK02899", "紫金矿业", 21.67, 1106.618293645365, 20, 22132.3658729073, 102146.0826454179),
    @("HK06881", "中国银河", 5.02, 610.1281790147427, 8.4, 5125.076703723839, 102146.0826454179),
    @("100000", "现金", 29.32, 29945.09102604357, 1, 29945.09102604357, 102146.0826454179)
)
Add-PortfolioRows $ws "大智 (稳健智远)" 41 $timestamp $rows

# Sheet 2: 大成投资组合收益 (新增第30-36行)
$ws = $wb.Worksheets.Item(2)
$rows = @(
    @("000725", "京东方A", 4.91, 1243.91705951017, 3.93, 4888.594043874969, 99563.82109128355),
    @("159781", "科创创业ETF", 4.86, 9122.058436407913, 0.53, 4834.690971296194, 99563.82109128355),
    @("513100", "纳指ETF", 4.95, 3137.523283860047, 1.57, 4925.911555660274, 99563.82109128355),
    @("513290", "纳指生物科技ETF", 0.99, 879.627063510763, 1.12, 985.1823111320547, 99563.82109128355),
    @("603119", "浙江荣泰", 45.85, 1051.546584462582, 43.41, 45647.63723152068, 99563.82109128355),
    @("688290", "景业智能", 8.76, 161.2702946560293, 54.11, 8726.335643837745, 99563.82109128355),
    @("100000", "现金", 29.68, 29555.46933396164, 1, 29555.46933396164, 99563.82109128355)
)
Add-PortfolioRows $ws "大成 (锐进先锋)" 30 $timestamp $rows

# Sheet 3: 我的投资组合收益 (新增第61-75行)
$ws = $wb.Worksheets.Item(3)
$rows = @(
    @("000333", "美的集团", 0.96, 13.2658076636599, 72.2, 957.7913133162448, 99986.84367663674),
    @("000725", "京东方A", 4.96, 1262.931689191865, 3.93, 4963.321538524029, 99986.84367663674),
    @("159781", "科创创业ETF", 4.91, 9261.499054073673, 0.53, 4908.594498659047, 99986.84367663674),
    @("510050", "上证50ETF", 5, 1812.03242362311, 2.76, 5001.209489199783, 99986.84367663674),
    @("510300", "沪深300ETF", 5.01, 1250.302372299946, 4.01, 5013.712512922783, 99986.84367663674),
    @("513100", "纳指ETF", 1, 637.096750216533, 1.57, 1000.241897839957, 99986.84367663674),
    @("513290", "纳指生物科技ETF", 1, 893.0731230713899, 1.12, 1000.241897839957, 99986.84367663674),
    @("518880", "黄金ETF", 1.01, 135.1678240324266, 7.5, 1013.7586802432, 99986.84367663674),
    @("600085", "同仁堂", 1, 26.53865475829018, 37.53, 995.9957130786304, 99986.84367663674),
    @("600900", "长江电力", 1.01, 33.36363902067901, 30.19, 1007.248262034299, 99986.84367663674),
    @("600989", "宝丰能源", 4.94, 307.3884135955614, 16.08, 4942.805690616627, 99986.84367663674),
    @("603119", "浙江荣泰", 1.03, 23.72490270018873, 43.41, 1029.898026215193, 99986.84367663674),
    @("HK02899", "紫金矿业", 1.11, 55.44578147671601, 20, 1108.91562953432, 99986.84367663674),
    @("HK06881", "中国银河", 1.03, 122.2789606161316, 8.4, 1027.143269175506, 99986.84367663674),
    @("100000", "现金", 66.02, 66015.96525743716, 1, 66015.96525743716, 99986.84367663674)
)
Add-PortfolioRows $ws "范式进化投资组合" 61 $timestamp $rows
